# CDIN.xlsx update:
#  - rename "Groupe" labels 2-A/2-B/2-C -> 1-A/1-B/1-C (all students use the
#    same shared label, so every row in that group is updated)
#  - bump the "Numero" (student id) column from 2015xxx -> 2017xxx
#  - refresh the "Moyenne de l'etudiant" column with the new values
#    ("correctif ... mise a jour semestre etudiant lors de l'inscription")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Groupe relabelling -----------------------------------------------
$group1A = @(5,6,7,16,23,24,27,31,35,37,38,39,41,43,44,46,52,57,60,62)
$group1B = @(3,8,11,15,19,20,21,25,26,29,30,33,40,47,49,50,53,56,58,59,63)
$group1C = @(4,9,10,12,13,14,17,18,22,28,32,34,36,42,45,48,51,54,55,61)

foreach ($r in $group1A) { $ws.Range("D$r").Value = "1-A" }
foreach ($r in $group1B) { $ws.Range("D$r").Value = "1-B" }
foreach ($r in $group1C) { $ws.Range("D$r").Value = "1-C" }

# --- Numero (A) + Moyenne de l'etudiant (E) ----------------------------
$rowUpdates = @(
    @{Row=3;  A=20170926; E=17},
    @{Row=4;  A=20170927; E=7},
    @{Row=5;  A=20170928; E=14},
    @{Row=6;  A=20170929; E=10},
    @{Row=7;  A=20170930; E=17},
    @{Row=8;  A=20170931; E=11},
    @{Row=9;  A=20170932; E=8},
    @{Row=10; A=20170933; E=8},
    @{Row=11; A=20170934; E=7},
    @{Row=12; A=20170935},
    @{Row=13; A=20170936; E=14},
    @{Row=14; A=20170937; E=5},
    @{Row=15; A=20170938; E=13},
    @{Row=16; A=20170939; E=20},
    @{Row=17; A=20170940; E=12},
    @{Row=18; A=20170941; E=14},
    @{Row=19; A=20170942; E=13},
    @{Row=20; A=20170943; E=13},
    @{Row=21; A=20170944; E=13},
    @{Row=22; A=20170945; E=20},
    @{Row=23; A=20170946; E=9},
    @{Row=24; A=20170947; E=7},
    @{Row=25; A=20170948; E=9},
    @{Row=26; A=20170949; E=11},
    @{Row=27; A=20170950; E=12},
    @{Row=28; A=20170951; E=15},
    @{Row=29; A=20170952; E=16},
    @{Row=30; A=20170953; E=18},
    @{Row=31; A=20170954; E=13},
    @{Row=32; A=20170955; E=12},
    @{Row=33; A=20170956; E=18},
    @{Row=34; A=20170957; E=16},
    @{Row=35; A=20170958},
    @{Row=36; A=20170959; E=17},
    @{Row=37; A=20170960; E=12},
    @{Row=38; A=20170961; E=15},
    @{Row=39; A=20170962; E=16},
    @{Row=40; A=20170963; E=16},
    @{Row=41; A=20170964; E=11},
    @{Row=42; A=20170965},
    @{Row=43; A=20170966; E=14},
    @{Row=44; A=20170967; E=10},
    @{Row=45; A=20170968; E=12},
    @{Row=46; A=20170969; E=13},
    @{Row=47; A=20170970; E=9},
    @{Row=48; A=20170971; E=8},
    @{Row=49; A=20170972; E=8},
    @{Row=50; A=20170973; E=5},
    @{Row=51; A=20170974; E=17},
    @{Row=52; A=20170975; E=5},
    @{Row=53; A=20170976; E=13},
    @{Row=54; A=20170977; E=6},
    @{Row=55; A=20170978; E=10},
    @{Row=56; A=20170979; E=14},
    @{Row=57; A=20170980; E=15},
    @{Row=58; A=20170981; E=8},
    @{Row=59; A=20170982; E=11},
    @{Row=60; A=20170983; E=17},
    @{Row=61; A=20170984; E=14},
    @{Row=62; A=20170985; E=6},
    @{Row=63; A=20170986}
)

foreach ($u in $rowUpdates) {
    $ws.Range("A$($u.Row)").Value = $u.A
    if ($u.ContainsKey("E")) {
        $ws.Range("E$($u.Row)").Value = $u.E
    }
}
